$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functies")

# Capture source values (from rows that stay put, 1-10) before the row insert shifts anything.
$valMedewerker      = $ws.Range("B14").Value2   # "Beheer > Medewerkers > Medewerker"
$valPenO            = $ws.Range("G10").Value2   # "P&O"
$valCUD             = $ws.Range("I10").Value2   # "CUD"
$valArbeidsmod      = $ws.Range("J10").Value2   # "Arbeidsmodaliteit"
$valInvariant       = $ws.Range("L10").Value2   # Invariant: tijdlijnArbeidsmodaliteiten... text

# Row 9 description height shrinks slightly.
$ws.Rows.Item(9).RowHeight = 86.4

# Insert a new row at position 11; everything from row 11 down shifts to row 12 down,
# inheriting the row-10 formatting the way Excel normally does on a row insert.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new "Beheren arbeidsmodaliteiten" function.
$ws.Range("B11").Value = $valMedewerker
$ws.Range("C11").Value = "Beheren arbeidsmodaliteiten"
$ws.Range("D11").Value = "het muteren van arbeidsmodaliteiten van een specifieke medewerker"
$ws.Range("G11").Value = $valPenO
$ws.Range("I11").Value = $valCUD
$ws.Range("J11").Value = $valArbeidsmod
$ws.Range("L11").Value = $valInvariant

$ws.Rows.Item(11).RowHeight = 43.2
